$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new price records were added (weekly Fruta/hortaliza update). Insert two
# blank rows at 205-206, pushing the existing rows 205+ down to 207+, then
# populate the two new rows with the new data.
$ws.Range("A205:A206").EntireRow.Insert()

# Row 205: new "Especial" quality record for Provincia de Curicó
$ws.Cells.Item(205, 1).Value = 6
$ws.Cells.Item(205, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(205, 3).Value = "Metropolitana"
$ws.Cells.Item(205, 4).Value = 44543
$ws.Cells.Item(205, 5).Value = 13
$ws.Cells.Item(205, 6).Value = "Fruta"
$ws.Cells.Item(205, 7).Value = 100101
$ws.Cells.Item(205, 8).Value = "Berries"
$ws.Cells.Item(205, 9).Value = 100101001
$ws.Cells.Item(205, 10).Value = "Arándano (blue)"
$ws.Cells.Item(205, 11).Value = "Sin especificar"
$ws.Cells.Item(205, 12).Value = "Especial"
$ws.Cells.Item(205, 13).Value = 200
$ws.Cells.Item(205, 14).Value = 5000
$ws.Cells.Item(205, 15).Value = 5000
$ws.Cells.Item(205, 16).Value = 5000
$ws.Cells.Item(205, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(205, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(205, 19).Value = 2500
$ws.Cells.Item(205, 20).Value = 2

# Row 206: new "Primera" quality record for Provincia de Curicó
$ws.Cells.Item(206, 1).Value = 6
$ws.Cells.Item(206, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(206, 3).Value = "Metropolitana"
$ws.Cells.Item(206, 4).Value = 44543
$ws.Cells.Item(206, 5).Value = 13
$ws.Cells.Item(206, 6).Value = "Fruta"
$ws.Cells.Item(206, 7).Value = 100101
$ws.Cells.Item(206, 8).Value = "Berries"
$ws.Cells.Item(206, 9).Value = 100101001
$ws.Cells.Item(206, 10).Value = "Arándano (blue)"
$ws.Cells.Item(206, 11).Value = "Sin especificar"
$ws.Cells.Item(206, 12).Value = "Primera"
$ws.Cells.Item(206, 13).Value = 500
$ws.Cells.Item(206, 14).Value = 4000
$ws.Cells.Item(206, 15).Value = 4000
$ws.Cells.Item(206, 16).Value = 4000
$ws.Cells.Item(206, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(206, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(206, 19).Value = 2000
$ws.Cells.Item(206, 20).Value = 2
